$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 2.45
$ws.Range("H2").Value = 2.7
$ws.Range("I2").Value = 3.4
$ws.Range("J2").Value = 3.5
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 4.75
$ws.Range("S2").Value = 1.75
$ws.Range("T2").Value = 2.05
$ws.Range("X2").Value = 10
$ws.Range("Y2").Value = 12
$ws.Range("Z2").Value = 26
$ws.Range("AC2").Value = 4.75
$ws.Range("AE2").Value = 23
$ws.Range("AH2").Value = 6.5
$ws.Range("AP2").Value = 41
$ws.Range("AR2").Value = 126
$ws.Range("AU2").Value = 11

# Row 3 updates
$ws.Range("G3").Value = 1.55
$ws.Range("I3").Value = 5
$ws.Range("L3").Value = 5
$ws.Range("U3").Value = 1.62
$ws.Range("V3").Value = 2.2
$ws.Range("X3").Value = 9
$ws.Range("AA3").Value = 12
$ws.Range("AD3").Value = 8.5
$ws.Range("AU3").Value = 7.5
$ws.Range("AW3").Value = 7
$ws.Range("BA3").Value = 81

$wb.Save()
